# Refresh the cryptocurrency price/volume snapshot (Price column D, Volume(1h)
# column E) for rows 2-51 of Sheet1, matching the GitHub Actions data refresh.
# Price values that look like plain numbers are written with a leading
# apostrophe (then restyled "Normal") so they stay text and keep their exact
# printed digits (e.g. trailing zeros) instead of being coerced to a float.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.181.23"
$ws.Range("E2").Value = "  -0.17%  "
$ws.Range("D3").Value = "1.907.90"
$ws.Range("E3").Value = "  +1.67%  "
$ws.Range("E4").Value = "  +0.23%  "
$ws.Range("D5").Formula = "'314.37"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.00%  "
$ws.Range("D6").Formula = "'1.002"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.24%  "
$ws.Range("D7").Formula = "'0.5054"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.28%  "
$ws.Range("D8").Formula = "'0.3930"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.12%  "
$ws.Range("D9").Formula = "'0.09315"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -5.78%  "
$ws.Range("E10").Value = "  -0.33%  "
$ws.Range("D11").Formula = "'41.86"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.45%  "
$ws.Range("D12").Formula = "'6.402"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.77%  "
$ws.Range("D13").Formula = "'20.86"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.62%  "
$ws.Range("D14").Value = "1.908.71"
$ws.Range("E14").Value = "  +2.04%  "
$ws.Range("E15").Value = "  -2.00%  "
$ws.Range("D16").Formula = "'1.001"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.17%  "
$ws.Range("E17").Value = "  -1.00%  "
$ws.Range("D18").Formula = "'92.53"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.55%  "
$ws.Range("D19").Formula = "'0.06612"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.35%  "
$ws.Range("D20").Formula = "'17.98"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.07%  "
$ws.Range("D21").Formula = "'1.001"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.25%  "
$ws.Range("D22").Formula = "'6.209"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.52%  "
$ws.Range("D23").Value = "28.257.60"
$ws.Range("E23").Value = "  -0.18%  "
$ws.Range("E24").Value = "  +1.01%  "
$ws.Range("D25").Formula = "'2.319"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.53%  "
$ws.Range("D26").Formula = "'2.607"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.78%  "
$ws.Range("D27").Value = "2.130.71"
$ws.Range("E27").Value = "  +2.02%  "
$ws.Range("D28").Formula = "'21.04"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.95%  "
$ws.Range("D29").Formula = "'157.94"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.38%  "
$ws.Range("D30").Formula = "'127.36"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.30%  "
$ws.Range("D31").Formula = "'1.102"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.92%  "
$ws.Range("E32").Value = "  +0.76%  "
$ws.Range("D33").Formula = "'5.639"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.09%  "
$ws.Range("D34").Formula = "'3.613"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.19%  "
$ws.Range("D35").Formula = "'9.669"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.39%  "
$ws.Range("D36").Formula = "'0.06651"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.26%  "
$ws.Range("D37").Formula = "'0.02423"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.54%  "
$ws.Range("D38").Formula = "'1.243"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.54%  "
$ws.Range("D39").Formula = "'0.2189"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.06%  "
$ws.Range("D40").Formula = "'1.283"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +8.15%  "
$ws.Range("D41").Formula = "'0.6431"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.77%  "
$ws.Range("D42").Formula = "'5.008"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.17%  "
$ws.Range("D43").Formula = "'11.49"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.41%  "
$ws.Range("E44").Value = "  +0.29%  "
$ws.Range("D45").Formula = "'13.38"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.61%  "
$ws.Range("D46").Formula = "'0.6017"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Formula = "'3.718"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.65%  "
$ws.Range("D48").Formula = "'1.280"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.19%  "
$ws.Range("D49").Formula = "'2.020"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.70%  "
$ws.Range("D50").Formula = "'122.97"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Formula = "'1.185"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.91%  "
